# Applies the leve-profit recalculation updates captured in the commit diff.
# Each block targets one (sheet, row) pair from the Table_<Job> listings and
# rewrites the currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N).

$wb = $excel.ActiveWorkbook

# ALC!row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 157.66667
$ws.Cells.Item(33, 9).Value = 157.66667
$ws.Cells.Item(33, 11).Value = 157.66667
$ws.Cells.Item(33, 13).Value = 71.33332999999999

# ALC!row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 4005
$ws.Cells.Item(116, 9).Value = 4005
$ws.Cells.Item(116, 11).Value = 4005
$ws.Cells.Item(116, 13).Value = -563

# ALC!row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(129, 8).Value = 3021.6667
$ws.Cells.Item(129, 9).Value = 721.2
$ws.Cells.Item(129, 11).Value = 2163.6
$ws.Cells.Item(129, 13).Value = 2836.4

# ALC!row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 2745.15
$ws.Cells.Item(137, 9).Value = 2031.3846
$ws.Cells.Item(137, 11).Value = 6094.1538
$ws.Cells.Item(137, 13).Value = -3544.1538

# ARM!row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3307
$ws.Cells.Item(61, 9).Value = 2994.4285
$ws.Cells.Item(61, 10).Value = 5495
$ws.Cells.Item(61, 11).Value = 2994.4285
$ws.Cells.Item(61, 12).Value = 5495
$ws.Cells.Item(61, 13).Value = -2782.4285
$ws.Cells.Item(61, 14).Value = -5919

# ARM!row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(97, 8).Value = 829.5454999999999
$ws.Cells.Item(97, 9).Value = 712.5
$ws.Cells.Item(97, 11).Value = 712.5
$ws.Cells.Item(97, 13).Value = -216.5

# ARM!row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 0
$ws.Cells.Item(102, 9).Value = 0
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 11).Value = 0
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 13).ClearContents()
$ws.Cells.Item(102, 14).ClearContents()

# ARM!row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 1618.3889
$ws.Cells.Item(122, 9).Value = 1181.0769
$ws.Cells.Item(122, 11).Value = 3543.2307
$ws.Cells.Item(122, 13).Value = -1093.2307

# ARM!row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 3307
$ws.Cells.Item(136, 9).Value = 2994.4285
$ws.Cells.Item(136, 10).Value = 5495
$ws.Cells.Item(136, 11).Value = 8983.2855
$ws.Cells.Item(136, 12).Value = 16485
$ws.Cells.Item(136, 13).Value = -6433.2855
$ws.Cells.Item(136, 14).Value = -21585

# BSM!row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 1527.8182
$ws.Cells.Item(99, 9).Value = 1463.625
$ws.Cells.Item(99, 10).Value = 1699
$ws.Cells.Item(99, 11).Value = 1463.625
$ws.Cells.Item(99, 12).Value = 1699
$ws.Cells.Item(99, 13).Value = 34.375
$ws.Cells.Item(99, 14).Value = -4695

# BSM!row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 0
$ws.Cells.Item(105, 9).Value = 0
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 11).Value = 0
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 13).ClearContents()
$ws.Cells.Item(105, 14).ClearContents()

# BSM!row 135
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(135, 8).Value = 135255
$ws.Cells.Item(135, 10).Value = 135255
$ws.Cells.Item(135, 12).Value = 135255
$ws.Cells.Item(135, 14).Value = -145395

# CRP!row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1273.25
$ws.Cells.Item(16, 9).Value = 1197.6666
$ws.Cells.Item(16, 11).Value = 1197.6666
$ws.Cells.Item(16, 13).Value = -910.6666

# CRP!row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6220.7417
$ws.Cells.Item(31, 9).Value = 2676
$ws.Cells.Item(31, 10).Value = 7254.625
$ws.Cells.Item(31, 11).Value = 2676
$ws.Cells.Item(31, 12).Value = 7254.625
$ws.Cells.Item(31, 13).Value = -2381
$ws.Cells.Item(31, 14).Value = -7844.625

# CRP!row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 6220.7417
$ws.Cells.Item(34, 9).Value = 2676
$ws.Cells.Item(34, 10).Value = 7254.625
$ws.Cells.Item(34, 11).Value = 2676
$ws.Cells.Item(34, 12).Value = 7254.625
$ws.Cells.Item(34, 13).Value = -2474
$ws.Cells.Item(34, 14).Value = -7658.625

# CRP!row 51
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(51, 8).Value = 74000
$ws.Cells.Item(51, 9).Value = 74000
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 11).Value = 74000
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(51, 13).Value = -73264
$ws.Cells.Item(51, 14).ClearContents()

# CRP!row 61
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(61, 8).Value = 74000
$ws.Cells.Item(61, 9).Value = 74000
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 74000
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).Value = -73652
$ws.Cells.Item(61, 14).ClearContents()

# CRP!row 87
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 12).Value = 0
$ws.Cells.Item(87, 14).ClearContents()

# CRP!row 90
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(90, 8).Value = 0
$ws.Cells.Item(90, 10).Value = 0
$ws.Cells.Item(90, 12).Value = 0
$ws.Cells.Item(90, 14).ClearContents()

# CRP!row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 2335.1667
$ws.Cells.Item(99, 9).Value = 1502.75
$ws.Cells.Item(99, 11).Value = 1502.75
$ws.Cells.Item(99, 13).Value = -4.75

# CRP!row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(113, 8).Value = 1273.25
$ws.Cells.Item(113, 9).Value = 1197.6666
$ws.Cells.Item(113, 11).Value = 1197.6666
$ws.Cells.Item(113, 13).Value = 972.3334

# CRP!row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 1730.1428
$ws.Cells.Item(122, 9).Value = 2299.6667
$ws.Cells.Item(122, 10).Value = 1303
$ws.Cells.Item(122, 11).Value = 6899.000100000001
$ws.Cells.Item(122, 12).Value = 3909
$ws.Cells.Item(122, 13).Value = -4449.000100000001
$ws.Cells.Item(122, 14).Value = -8809

# CRP!row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 2335.1667
$ws.Cells.Item(126, 9).Value = 1502.75
$ws.Cells.Item(126, 11).Value = 4508.25
$ws.Cells.Item(126, 13).Value = -2038.25

# CUL!row 14
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 470.30768
$ws.Cells.Item(14, 9).Value = 470.30768
$ws.Cells.Item(14, 11).Value = 1410.92304
$ws.Cells.Item(14, 13).Value = -1237.92304

# GSM!row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 4099
$ws.Cells.Item(80, 9).Value = 4099
$ws.Cells.Item(80, 11).Value = 4099
$ws.Cells.Item(80, 13).Value = -3101

# GSM!row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 4099
$ws.Cells.Item(83, 9).Value = 4099
$ws.Cells.Item(83, 11).Value = 20495
$ws.Cells.Item(83, 13).Value = -15503

# GSM!row 92
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(92, 8).Value = 7604.8
$ws.Cells.Item(92, 10).Value = 7604.8
$ws.Cells.Item(92, 12).Value = 7604.8
$ws.Cells.Item(92, 14).Value = -11348.8

# GSM!row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1267.3572
$ws.Cells.Item(102, 9).Value = 1326.4615
$ws.Cells.Item(102, 11).Value = 1326.4615
$ws.Cells.Item(102, 13).Value = 295.5385000000001

# GSM!row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 3467.8462
$ws.Cells.Item(107, 9).Value = 2598
$ws.Cells.Item(107, 11).Value = 2598
$ws.Cells.Item(107, 13).Value = -678

# LTW!row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 938.2727
$ws.Cells.Item(22, 9).Value = 904.6
$ws.Cells.Item(22, 11).Value = 904.6
$ws.Cells.Item(22, 13).Value = -609.6

# LTW!row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 938.2727
$ws.Cells.Item(27, 9).Value = 904.6
$ws.Cells.Item(27, 11).Value = 904.6
$ws.Cells.Item(27, 13).Value = -797.6

# LTW!row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 4887
$ws.Cells.Item(40, 9).Value = 4887
$ws.Cells.Item(40, 11).Value = 4887
$ws.Cells.Item(40, 13).Value = -4751

# LTW!row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 3249.75
$ws.Cells.Item(46, 10).Value = 5550
$ws.Cells.Item(46, 12).Value = 5550
$ws.Cells.Item(46, 14).Value = -5926

# LTW!row 94
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(94, 8).Value = 21000
$ws.Cells.Item(94, 10).Value = 21000
$ws.Cells.Item(94, 12).Value = 21000
$ws.Cells.Item(94, 14).Value = -22352

# LTW!row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 0
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 13).ClearContents()

# LTW!row 133
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(133, 8).Value = 100000
$ws.Cells.Item(133, 10).Value = 100000
$ws.Cells.Item(133, 12).Value = 100000
$ws.Cells.Item(133, 14).Value = -105060

# WVR!row 119
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(119, 8).Value = 49899.332
$ws.Cells.Item(119, 10).Value = 49899.332
$ws.Cells.Item(119, 12).Value = 49899.332
$ws.Cells.Item(119, 14).Value = -59575.332

# WVR!row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 3260.3845
$ws.Cells.Item(136, 9).Value = 2388
$ws.Cells.Item(136, 11).Value = 7164
$ws.Cells.Item(136, 13).Value = -4614

